$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Planilha1")
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "Planilha2"

# --- Row labels / shared strings, set in the same order as the source workbook ---
$ws2.Range("C4").Value = "Umidade máxima"
$ws2.Range("C8").Value = "Umidade mínima"
$ws2.Range("C5").Value = "Alerta (muito umido)"
$ws2.Range("C7").Value = "Alerta (muito seco)"
$ws2.Range("L4").Value = "Umidade:"
$ws2.Range("N4").Value = "muito seco"
$ws2.Range("P4").Value = "muito umido"
$ws2.Range("C6").Value = "Dado (temperatura)"

# --- Humidity data table C4:J8 ---
$ws2.Range("D4:J4").Value = 0.8
$ws2.Range("D5:J5").Value = 0.75
$ws2.Range("D6").Value = 0.71
$ws2.Range("E6").Value = 0.65
$ws2.Range("F6").Value = 0.62
$ws2.Range("G6").Value = 0.66
$ws2.Range("H6").Value = 0.76
$ws2.Range("I6").Value = 0.74
$ws2.Range("J6").Value = 0.73
$ws2.Range("D7:J7").Value = 0.65
$ws2.Range("D8:J8").Value = 0.6
$ws2.Range("D4:J8").NumberFormat = "0%"

# --- Legend table L4:Q5 (mirrors the "holder" legend on Planilha1) ---
$ws2.Range("M4").Value = "abaixo do limite"
$ws2.Range("O4").Value = "ideal"
$ws2.Range("Q4").Value = "acima do limite"

$ws2.Range("M4").Interior.Color = 0xA03070
$ws2.Range("M4").HorizontalAlignment = -4108
$ws2.Range("N4").Interior.Color = 0xF0B000
$ws2.Range("N4").HorizontalAlignment = -4108
$ws2.Range("O4").Interior.Color = 0x50D092
$ws2.Range("O4").HorizontalAlignment = -4108
$ws2.Range("P4").Interior.Color = 0x00C0FF
$ws2.Range("P4").HorizontalAlignment = -4108
$ws2.Range("Q4").Interior.Color = 0x0000FF
$ws2.Range("Q4").HorizontalAlignment = -4108

$ws2.Range("M5").Value = 0.6
$ws2.Range("N5").Value = 0.65
$ws2.Range("O5").Value = 0.7
$ws2.Range("P5").Value = 0.75
$ws2.Range("Q5").Value = 0.8

$ws2.Range("M5").Interior.Color = 0xA03070
$ws2.Range("M5").HorizontalAlignment = -4108
$ws2.Range("N5").Interior.Color = 0xF0B000
$ws2.Range("N5").HorizontalAlignment = -4108
$ws2.Range("O5").Interior.Color = 0x50D092
$ws2.Range("O5").HorizontalAlignment = -4108
$ws2.Range("P5").Interior.Color = 0x00C0FF
$ws2.Range("P5").HorizontalAlignment = -4108
$ws2.Range("Q5").Interior.Color = 0x0000FF
$ws2.Range("Q5").HorizontalAlignment = -4108
$ws2.Range("M5:Q5").NumberFormat = "0%"

Write-Host "done"
